$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column width adjustments ---
# Column B loses its "best fit" auto-width and gets a slightly wider explicit width.
$ws.Columns.Item(2).ColumnWidth = 9.3
# Column H gets wider to fit the new formula text.
$ws.Columns.Item(8).ColumnWidth = 34.3
# New columns I and J get explicit widths (used by the new SelectFirst tables).
$ws.Columns.Item(9).ColumnWidth = 11.5
$ws.Columns.Item(10).ColumnWidth = 13.2

# --- New "SelectFirst" Spreadsheet test tables (rows 27-30) ---

# Row 27: section headers, merged across B:D and G:I, left aligned.
$ws.Range("B27:D27").Value = "Spreadsheet Driver checkSpreadsheet1(Driver[] dd, int maxAge)"
$ws.Range("B27:D27").HorizontalAlignment = -4131
$ws.Range("B27:D27").Merge()

$ws.Range("G27:I27").Value = "Spreadsheet Driver checkSpreadsheet2(Driver[] dd, int dIndex)"
$ws.Range("G27:I27").HorizontalAlignment = -4131
$ws.Range("G27:I27").Merge()

# Row 28: spreadsheet column headers.
$ws.Range("C28").Value = "Formula:Driver"
$ws.Range("D28").Value = "Age:int"
$ws.Range("H28").Value = "Formula:Driver"
$ws.Range("I28").Value = "Index:int"

# Row 29: SelectStep row with formula-like text (entered as literal text, quote-prefixed).
$ws.Range("B29").Value = "SelectStep"
$ws.Range("C29").Value = "'=dd[!@ age < `$Age]"
$ws.Range("D29").Value = "'=maxAge"

$ws.Range("G29").Value = "SelectStep"
$ws.Range("H29").Value = "'=dd[!@ name == testDrivers[`$Index].name]"
$ws.Range("I29").Value = "'=dIndex"

# Row 30: RETURN row.
$ws.Range("B30").Value = "RETURN"
$ws.Range("C30").Value = "'=`$SelectStep"

$ws.Range("G30").Value = "RETURN"
$ws.Range("H30").Value = "'=`$SelectStep"

# --- Selection mirrors where the author ended up after adding the tables ---
$ws.Range("G27:I27").Select()
